$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 04:52"

# --- Australia (row 47) gets refreshed totals ---
$ws.Range("B47").Value = 6752
$ws.Range("C47").Value = 6
$ws.Range("D47").Value = 5715
$ws.Range("E47").Value = 946
$ws.Range("F47").Value = 38
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 91

# --- New country "Guatemala" inserted between Uruguay and Somalia ---
# Somalia & San Marino shift down one row (106->107, 107->108)
$ws.Range("A108").Value = "San Marino"
$ws.Range("B108").Value = 563
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 69
$ws.Range("E108").Value = 453
$ws.Range("F108").Value = 6
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 41

$ws.Range("A107").Value = "Somalia"
$ws.Range("B107").Value = 582
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 20
$ws.Range("E107").Value = 534
$ws.Range("F107").Value = 2
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 28

$ws.Range("A106").Value = "Guatemala"
$ws.Range("B106").Value = 585
$ws.Range("C106").Value = 28
$ws.Range("D106").Value = 65
$ws.Range("E106").Value = 504
$ws.Range("F106").Value = 5
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 16

# --- New country "Yemen" inserted between Montserrat and Burundi ---
# Burundi..Butan shift down one row (200->201 ... 209->210)
$ws.Range("A210").Value = "Butan"
$ws.Range("B210").Value = 7
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 5
$ws.Range("E210").Value = 2
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

$ws.Range("A209").Value = "Mauritania"
$ws.Range("B209").Value = 8
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 6
$ws.Range("E209").Value = 1
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1

$ws.Range("A208").Value = "Santo Tome y Principe"
$ws.Range("B208").Value = 8
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 4
$ws.Range("E208").Value = 4
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

$ws.Range("A207").Value = "Papua Nueva Guinea"
$ws.Range("B207").Value = 8
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 0
$ws.Range("E207").Value = 8
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

$ws.Range("A206").Value = "Gambia"
$ws.Range("B206").Value = 10
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 8
$ws.Range("E206").Value = 1
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 1

$ws.Range("A205").Value = "Surinam"
$ws.Range("B205").Value = 10
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 8
$ws.Range("E205").Value = 1
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 1

$ws.Range("A204").Value = "Santa Sede"
$ws.Range("B204").Value = 10
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 2
$ws.Range("E204").Value = 8
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

$ws.Range("A203").Value = "Groenlandia"
$ws.Range("B203").Value = 11
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 11
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

$ws.Range("A202").Value = "Seychelles"
$ws.Range("B202").Value = 11
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 6
$ws.Range("E202").Value = 5
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A201").Value = "Burundi"
$ws.Range("B201").Value = 11
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 4
$ws.Range("E201").Value = 6
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 1

$ws.Range("A200").Value = "Yemen"
$ws.Range("B200").Value = 11
$ws.Range("C200").Value = 5
$ws.Range("D200").Value = 1
$ws.Range("E200").Value = 8
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 2
$ws.Range("H200").Value = 2
